$wb = $excel.ActiveWorkbook

# --- Update selection on the GCMS-Medusa sheet (it will lose tabSelected
#     automatically once a different sheet becomes active) ---
$medusa = $wb.Worksheets.Item("GCMS-Medusa")
$medusa.Range("B7").Select() | Out-Null

# --- Add the new "Picarro" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Picarro"

# Reference sheets whose existing cell formats/styles we reuse
$ale = $wb.Worksheets.Item("ALE")

# --- Copy cell formatting (styles) from existing sheets that already use
#     the exact same style combinations required for the new sheet ---
$ale.Range("A1:A6").Copy()
$newSheet.Range("A1:A6").PasteSpecial(-4122)

$ale.Range("A7").Copy()
$newSheet.Range("A7").PasteSpecial(-4122)

$medusa.Range("B7").Copy()
$newSheet.Range("B7").PasteSpecial(-4122)

$ale.Range("A8").Copy()
$newSheet.Range("A8").PasteSpecial(-4122)

$ale.Range("B8").Copy()
$newSheet.Range("B8").PasteSpecial(-4122)

$ale.Range("A9").Copy()
$newSheet.Range("A9").PasteSpecial(-4122)

# --- Fill in the values for the new "Picarro" sheet ---
$newSheet.Range("A1").Value = "# AGAGE GCMD data release schedule"
$newSheet.Range("A2").Value = "# DO NOT CHANGE THE FORMAT OF THIS SPREADSHEET"
$newSheet.Range("A3").Value = "# ALL GRID CELLS MUST BE IN TEXT FORMAT (NOT DATE FORMAT)"
$newSheet.Range("A4").Value = "# Date format YYYY-MM-DD HH:MM"
$newSheet.Range("A5").Value = '# All data will be processed until the "general release date" unless specified in the table'
$newSheet.Range("A6").Value = '# "x" indicates that no part of the record will be processed'
$newSheet.Range("A7").Value = "General release date"
$newSheet.Range("B7").Value = "2023-01-01 00:00"
$newSheet.Range("A8").Value = "Species"
$newSheet.Range("B8").Value = "THD"
$newSheet.Range("A9").Value = "ch4"

# --- Column widths matching the other sheets (A=20, B=15.5 best-fit) ---
$newSheet.Columns.Item(1).ColumnWidth = 19.1666666666667
$newSheet.Columns.Item(2).ColumnWidth = 14.6666666666667

# --- Match final selection on the new sheet ---
$newSheet.Range("E9").Select() | Out-Null
